$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix typo in existing label "Comparision" -> "Comparison"
$ws.Range("C38").Value = "WebApp Comparison Overview"

# Add a new row of tracked time for the "Comparison Detail" work
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A39").Value = 45289
$ws.Range("B39").Value = 8
$ws.Range("C39").Value = "WebApp Comparison Detail"

# Match the selection state recorded after the edit
$ws.Range("C39").Select()
